$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-29 18:31:47"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
